$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the current data block (rows 74-75),
# shifting the existing rows 74-110 down to 76-112.
$ws.Range("A74:A75").EntireRow.Insert()

# Row 74 - new weekly record (Black Amber / Primera)
$ws.Cells.Item(74, 1).Value = 4
$ws.Cells.Item(74, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(74, 3).Value = "Los Lagos"
$ws.Cells.Item(74, 4).Value = 44572
$ws.Cells.Item(74, 5).Value = 10
$ws.Cells.Item(74, 6).Value = "Fruta"
$ws.Cells.Item(74, 7).Value = 100103
$ws.Cells.Item(74, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(74, 9).Value = 100103002
$ws.Cells.Item(74, 10).Value = "Ciruela"
$ws.Cells.Item(74, 11).Value = "Black Amber"
$ws.Cells.Item(74, 12).Value = "Primera"
$ws.Cells.Item(74, 13).Value = 600
$ws.Cells.Item(74, 14).Value = 18000
$ws.Cells.Item(74, 15).Value = 18500
$ws.Cells.Item(74, 16).Value = 18250
$ws.Cells.Item(74, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(74, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(74, 19).Value = 1217
$ws.Cells.Item(74, 20).Value = 15

# Row 75 - new weekly record (Black Amber / Segunda)
$ws.Cells.Item(75, 1).Value = 4
$ws.Cells.Item(75, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(75, 3).Value = "Los Lagos"
$ws.Cells.Item(75, 4).Value = 44572
$ws.Cells.Item(75, 5).Value = 10
$ws.Cells.Item(75, 6).Value = "Fruta"
$ws.Cells.Item(75, 7).Value = 100103
$ws.Cells.Item(75, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(75, 9).Value = 100103002
$ws.Cells.Item(75, 10).Value = "Ciruela"
$ws.Cells.Item(75, 11).Value = "Black Amber"
$ws.Cells.Item(75, 12).Value = "Segunda"
$ws.Cells.Item(75, 13).Value = 300
$ws.Cells.Item(75, 14).Value = 16000
$ws.Cells.Item(75, 15).Value = 16000
$ws.Cells.Item(75, 16).Value = 16000
$ws.Cells.Item(75, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(75, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(75, 19).Value = 1067
$ws.Cells.Item(75, 20).Value = 15
